$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to Text format for the new rows so the date-like strings
# ("YYYY-MM-DD") are stored as text, matching the existing rows in the sheet
# (which use inline/shared strings rather than date serial numbers).
$ws.Range("A461:A472").NumberFormat = "@"

$ws.Range("A461").Value = "2021-11-15"
$ws.Range("B461").Value = "overview"
$ws.Range("C461").Value = "K02000001"
$ws.Range("D461").Value = "United Kingdom"
$ws.Range("E461").Value = 9600369
$ws.Range("F461").Value = 39705
$ws.Range("G461").Value = 47
$ws.Range("H461").Value = 142945

$ws.Range("A462").Value = "2021-11-16"
$ws.Range("B462").Value = "overview"
$ws.Range("C462").Value = "K02000001"
$ws.Range("D462").Value = "United Kingdom"
$ws.Range("E462").Value = 9637190
$ws.Range("F462").Value = 37243
$ws.Range("G462").Value = 214
$ws.Range("H462").Value = 143159

$ws.Range("A463").Value = "2021-11-17"
$ws.Range("B463").Value = "overview"
$ws.Range("C463").Value = "K02000001"
$ws.Range("D463").Value = "United Kingdom"
$ws.Range("E463").Value = 9675058
$ws.Range("F463").Value = 38263
$ws.Range("G463").Value = 201
$ws.Range("H463").Value = 143360

$ws.Range("A464").Value = "2021-11-18"
$ws.Range("B464").Value = "overview"
$ws.Range("C464").Value = "K02000001"
$ws.Range("D464").Value = "United Kingdom"
$ws.Range("E464").Value = 9721916
$ws.Range("F464").Value = 46807
$ws.Range("G464").Value = 199
$ws.Range("H464").Value = 143559

$ws.Range("A465").Value = "2021-11-19"
$ws.Range("B465").Value = "overview"
$ws.Range("C465").Value = "K02000001"
$ws.Range("D465").Value = "United Kingdom"
$ws.Range("E465").Value = 9766153
$ws.Range("F465").Value = 44242
$ws.Range("G465").Value = 157
$ws.Range("H465").Value = 143716

$ws.Range("A466").Value = "2021-11-20"
$ws.Range("B466").Value = "overview"
$ws.Range("C466").Value = "K02000001"
$ws.Range("D466").Value = "United Kingdom"
$ws.Range("E466").Value = 9806034
$ws.Range("F466").Value = 40941
$ws.Range("G466").Value = 150
$ws.Range("H466").Value = 143866

$ws.Range("A467").Value = "2021-11-21"
$ws.Range("B467").Value = "overview"
$ws.Range("C467").Value = "K02000001"
$ws.Range("D467").Value = "United Kingdom"
$ws.Range("E467").Value = 9845492
$ws.Range("F467").Value = 40004
$ws.Range("G467").Value = 61
$ws.Range("H467").Value = 143927

$ws.Range("A468").Value = "2021-11-22"
$ws.Range("B468").Value = "overview"
$ws.Range("C468").Value = "K02000001"
$ws.Range("D468").Value = "United Kingdom"
$ws.Range("E468").Value = 9889926
$ws.Range("F468").Value = 44917
$ws.Range("G468").Value = 45
$ws.Range("H468").Value = 143972

$ws.Range("A469").Value = "2021-11-23"
$ws.Range("B469").Value = "overview"
$ws.Range("C469").Value = "K02000001"
$ws.Range("D469").Value = "United Kingdom"
$ws.Range("E469").Value = 9932408
$ws.Range("F469").Value = 42484
$ws.Range("G469").Value = 165
$ws.Range("H469").Value = 144137

$ws.Range("A470").Value = "2021-11-24"
$ws.Range("B470").Value = "overview"
$ws.Range("C470").Value = "K02000001"
$ws.Range("D470").Value = "United Kingdom"
$ws.Range("E470").Value = 9974843
$ws.Range("F470").Value = 43676
$ws.Range("G470").Value = 149
$ws.Range("H470").Value = 144286

$ws.Range("A471").Value = "2021-11-25"
$ws.Range("B471").Value = "overview"
$ws.Range("C471").Value = "K02000001"
$ws.Range("D471").Value = "United Kingdom"
$ws.Range("E471").Value = 10021497
$ws.Range("F471").Value = 47240
$ws.Range("G471").Value = 147
$ws.Range("H471").Value = 144433

$ws.Range("A472").Value = "2021-11-26"
$ws.Range("B472").Value = "overview"
$ws.Range("C472").Value = "K02000001"
$ws.Range("D472").Value = "United Kingdom"
$ws.Range("E472").Value = 10070841
$ws.Range("F472").Value = 50091
$ws.Range("G472").Value = 160
$ws.Range("H472").Value = 144593

Write-Output "Added rows 461-472 to covid_totals sheet"
